# costs_op.xlsx: inpatient and outpatient costs S1, set max duration IO therapy 2L+ to 2 years
#
# Row 2 (state S1): recompute the per-cycle mean cost (B2) and its recomputed
# standard-error-derived value (C2) using the new 2-year-duration inputs, and
# repoint the reference (D2) at the new BIA citation.
# Row 3 (state P1/S2): updated mean cost (B3) -> recalculated se (C3).
# Row 4: no value change, only picks up the same "Comma" number format as C2/C3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Number format string used by the existing "Comma"-styled cells (C3/C4);
# re-applying the exact formatCode (incl. the escaped parentheses) makes the
# engine reuse that existing cell format instead of registering a new one.
$commaFormat = '_(* #,##0_);_(* \(#,##0\);_(* "-"??_);_(@_)'

# --- Row 2 : S1 ---
$ws.Range("B2").Formula = "=ROUND(1.2*(1643.16-(1375.84+121.61+93.2)),0)"
$ws.Range("C2").Formula = "=ROUND(B2*((2927.15-1027.33)/3.92)/1971.79,0)"
$ws.Range("C2").NumberFormat = $commaFormat
$ws.Range("D2").Value = "Graham2018BIA"

# --- Row 3 : P1/S2 ---
$ws.Range("B3").Formula = "=1442+182"

# --- Row 4 : align number format with C2/C3 (value/formula unchanged) ---
$ws.Range("C4").NumberFormat = $commaFormat

# Restore the selection to match the saved workbook view (C5).
$ws.Range("C5").Select()
